$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact" row -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second duplicate "Contact" row -> delete it entirely (rows below shift up)
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

# K2 (Short): "Extension" -> "Claim Response Item Status"
$elem.Range("K2").Value = "Claim Response Item Status"

# L2 (Definition): "An Extension" -> "Payment status of claim item"
$elem.Range("L2").Value = "Payment status of claim item"
